$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (row index -> new text)
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "10247"
    5  = "0.00001"
    7  = "0.00017"
    8  = "0.00055"
    9  = "0.00028"
    10 = "0.00033"
    11 = "0.00042"
    12 = "1.87238"
}

foreach ($row in $updates.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $updates[$row]
}

# Rows 44-46 previously held multiple tab-separated values in one run;
# collapse each down to a single value.
$collapsed = @{
    44 = "99.78"
    45 = "1.87"
    46 = "844"
}

foreach ($row in $collapsed.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $collapsed[$row]
}
